# Commit: "Fruta / hortaliza, semanal"
# A new weekly price-report row is inserted at row 59 (pushing the existing
# rows 59-190 down to 60-191); the new row carries a fresh price quote for
# the same market/product series.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 59:190 down to 60:191, leaving row 59 free for the new record.
$ws.Rows.Item(59).Insert()

$ws.Range("A59").Value2 = 4
$ws.Range("B59").Value2 = "Feria Lagunitas de Puerto Montt"
$ws.Range("C59").Value2 = "Los Lagos"
$ws.Range("D59").Value2 = 44519
$ws.Range("E59").Value2 = 10
$ws.Range("F59").Value2 = 100112037
$ws.Range("G59").Value2 = "Cebollín"
$ws.Range("H59").Value2 = "Sin especificar"
$ws.Range("I59").Value2 = "Primera"
$ws.Range("J59").Value2 = 180
$ws.Range("K59").Value2 = 5000
$ws.Range("L59").Value2 = 5000
$ws.Range("M59").Value2 = 5000
$ws.Range("N59").Value2 = "$/paquete 36 unidades"
$ws.Range("O59").Value2 = "Región Metropolitana"
$ws.Range("P59").Value2 = 139
$ws.Range("Q59").Value2 = 36
$ws.Range("R59").Value2 = "Hortaliza"
